# feat: add 2022-Q3 data
#
# Inserts a brand-new "2022-Q3" sheet (duplicated from the existing
# "2022-Q2" sheet so it inherits the same column layout / styles), fills it
# with the 2022-Q3 fund-holdings data, and updates the "总计" summary sheet
# with a new leading row for 2022-Q3 (shifting the older rows down).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a value into a cell while forcing TEXT storage (so values
# like "89.00" or "013082" are not silently coerced into numbers and lose
# their formatting / leading zeros), without leaving a lingering custom
# number-format behind on the cell.
# ---------------------------------------------------------------------
function Set-TextCell($cell, [string]$val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------
# 1. Duplicate the "2022-Q2" sheet (current 2nd sheet) and place the copy
#    right before it; rename the copy to "2022-Q3". This keeps the same
#    header row / column styling that every quarterly sheet shares.
# ---------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item(2)
$q2Sheet.Copy($q2Sheet, $null)

$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# ---------------------------------------------------------------------
# 2. Populate the new "2022-Q3" sheet with this quarter's fund data.
#    (Columns: A idx, B code, C name, D size, E stock-position,
#     F position-share, G market-value, H position-rank)
# ---------------------------------------------------------------------
$q3Data = @(
    @(0,  "516970", "广发中证基建工程ETF",               "89.00", "99.56", "5.45", "4.8505", 6),
    @(1,  "165525", "信诚中证基建工程指数（LOF）",         "10.82", "94.26", "5.18", "0.5605", 6),
    @(2,  "516950", "银华中证基建ETF",                   "11.07", "97.93", "3.93", "0.4351", 7),
    @(3,  "013082", "信诚中证基建工程指数（LOF）C",        "3.48",  "94.26", "5.18", "0.1803", 6),
    @(4,  "159635", "华夏中证基建ETF",                   "3.40",  "99.03", "3.98", "0.1353", 7),
    @(5,  "159619", "国泰中证基建ETF",                   "3.30",  "98.76", "3.97", "0.1310", 6),
    @(6,  "013326", "万家景气驱动混合A",                 "2.03",  "92.93", "3.17", "0.0644", 8),
    @(7,  "005933", "新疆前海联合先进制造灵活配置混合A",   "0.95",  "92.14", "5.12", "0.0486", 3),
    @(8,  "013757", "泰信均衡价值混合A",                 "0.76",  "66.30", "3.18", "0.0242", 10),
    @(9,  "008491", "万家周期优势企业混合A",             "0.61",  "93.50", "2.62", "0.0160", 10),
    @(10, "013327", "万家景气驱动混合C",                 "0.32",  "92.93", "3.17", "0.0101", 8),
    @(11, "013758", "泰信均衡价值混合C",                 "0.30",  "66.30", "3.18", "0.0095", 10),
    @(12, "005934", "新疆前海联合先进制造灵活配置混合C",   "0.10",  "92.14", "5.12", "0.0051", 3),
    @(13, "008492", "万家周期优势企业混合C",             "0.14",  "93.50", "2.62", "0.0037", 10)
)

$lastExistingRow = $q3Sheet.UsedRange.Rows.Count   # 13 (header + 12 rows), inherited from 2022-Q2

for ($i = 0; $i -lt $q3Data.Count; $i++) {
    $row = $i + 2
    $vals = $q3Data[$i]

    if ($row -gt $lastExistingRow) {
        # Need a brand-new row: copy formatting (esp. the "A" index-column
        # style) down from the row above so it matches the sheet's pattern.
        $q3Sheet.Cells.Item($row - 1, 1).Copy($q3Sheet.Cells.Item($row, 1))
    }

    $q3Sheet.Cells.Item($row, 1).Value = $vals[0]          # A: index (number)
    Set-TextCell $q3Sheet.Cells.Item($row, 2) $vals[1]      # B: fund code (text)
    Set-TextCell $q3Sheet.Cells.Item($row, 3) $vals[2]      # C: fund name (text)
    Set-TextCell $q3Sheet.Cells.Item($row, 4) $vals[3]      # D: fund size (text)
    Set-TextCell $q3Sheet.Cells.Item($row, 5) $vals[4]      # E: stock position (text)
    Set-TextCell $q3Sheet.Cells.Item($row, 6) $vals[5]      # F: position share (text)
    Set-TextCell $q3Sheet.Cells.Item($row, 7) $vals[6]      # G: market value (text)
    $q3Sheet.Cells.Item($row, 8).Value = $vals[7]           # H: position rank (number)
}

# ---------------------------------------------------------------------
# 3. Update the "总计" (summary) sheet: shift the 4 existing rows down by
#    one and insert the new 2022-Q3 totals in row 2.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)

# Copy row 5 -> row 6 (keeps the "A" column style identical to the source
# cell instead of leaving a brand-new unstyled cell behind).
$totalSheet.Cells.Item(5, 1).Copy($totalSheet.Cells.Item(6, 1))
$totalSheet.Cells.Item(6, 1).Value = 4
$totalSheet.Cells.Item(6, 2).Value = "2021-Q3"
$totalSheet.Cells.Item(6, 3).Value = 12
$totalSheet.Cells.Item(6, 4).Value = 10.16

$totalSheet.Cells.Item(5, 1).Value = 3
$totalSheet.Cells.Item(5, 2).Value = "2021-Q4"
$totalSheet.Cells.Item(5, 3).Value = 7
$totalSheet.Cells.Item(5, 4).Value = 1.36

$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(4, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(4, 3).Value = 6
$totalSheet.Cells.Item(4, 4).Value = 4.69

$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(3, 2).Value = "2022-Q2"
$totalSheet.Cells.Item(3, 3).Value = 12
$totalSheet.Cells.Item(3, 4).Value = 7.76

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 14
$totalSheet.Cells.Item(2, 4).Value = 6.47

# ---------------------------------------------------------------------
# 4. Restore the originally active sheet ("2021-Q3", now the last tab)
#    since inserting/copying sheets shifts the active-tab focus.
# ---------------------------------------------------------------------
$wb.Worksheets.Item($wb.Worksheets.Count).Activate()
